# Auto-generated edit script: applies the 2022-07-19 daily crime data update
# to the "violent-crime-full-year" workbook across all affected sheets.

$wb = $excel.ActiveWorkbook

# Sheet 1: Citywide Totals
$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('I2').Value = 3775
$ws.Range('I3').Value = 3903
$ws.Range('G4').Value = 1437
$ws.Range('H4').Value = 1667
$ws.Range('I4').Value = 918
$ws.Range('I5').Value = 361
$ws.Range('I6').Value = 4383
$ws.Range('G7').Value = 24660
$ws.Range('H7').Value = 25977
$ws.Range('I7').Value = 13340

# Sheet 2: By Neighborhood
$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('I4').Value = 51
$ws.Range('I6').Value = 91
$ws.Range('I8').Value = 804
$ws.Range('I12').Value = 27
$ws.Range('I14').Value = 68
$ws.Range('I19').Value = 360
$ws.Range('I20').Value = 326
$ws.Range('I21').Value = 74
$ws.Range('I23').Value = 128
$ws.Range('I29').Value = 863
$ws.Range('I30').Value = 47
$ws.Range('I33').Value = 601
$ws.Range('I36').Value = 185
$ws.Range('I37').Value = 433
$ws.Range('I41').Value = 59
$ws.Range('I42').Value = 464
$ws.Range('I43').Value = 115
$ws.Range('I44').Value = 96
$ws.Range('I47').Value = 94
$ws.Range('I49').Value = 110
$ws.Range('I54').Value = 305
$ws.Range('I60').Value = 66
$ws.Range('G63').Value = 201
$ws.Range('H63').Value = 208
$ws.Range('I63').Value = 46
$ws.Range('I65').Value = 296
$ws.Range('I66').Value = 36
$ws.Range('I67').Value = 512
$ws.Range('I72').Value = 49
$ws.Range('I73').Value = 111
$ws.Range('I76').Value = 202
$ws.Range('I78').Value = 193
$ws.Range('I79').Value = 361
$ws.Range('I83').Value = 271
$ws.Range('I85').Value = 602
$ws.Range('I86').Value = 84
$ws.Range('I89').Value = 149
$ws.Range('I91').Value = 164
$ws.Range('I92').Value = 41
$ws.Range('I94').Value = 127
$ws.Range('I95').Value = 213
$ws.Range('I96').Value = 148
$ws.Range('I97').Value = 101
$ws.Range('I98').Value = 88
$ws.Range('I99').Value = 248
$ws.Range('G101').Value = 24660
$ws.Range('H101').Value = 25977
$ws.Range('I101').Value = 13340

# Sheet 3: South Shore
$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('I3').Value = 242
$ws.Range('I4').Value = 36
$ws.Range('I5').Value = 20
$ws.Range('I6').Value = 148
$ws.Range('I7').Value = 602

# Sheet 7: Austin
$ws = $wb.Worksheets.Item('Austin')
$ws.Range('I2').Value = 250
$ws.Range('I3').Value = 223
$ws.Range('I7').Value = 804

# Sheet 10: Uptown
$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('I6').Value = 55
$ws.Range('I7').Value = 149

# Sheet 11: West Ridge
$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range('I2').Value = 41
$ws.Range('I4').Value = 8
$ws.Range('I7').Value = 148

# Sheet 12: Bridgeport
$ws = $wb.Worksheets.Item('Bridgeport')
$ws.Range('I3').Value = 19
$ws.Range('I7').Value = 68

# Sheet 13: Fuller Park
$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('I2').Value = 11
$ws.Range('I3').Value = 16
$ws.Range('I7').Value = 47

# Sheet 14: Grand Crossing
$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('I2').Value = 137
$ws.Range('I3').Value = 137
$ws.Range('I7').Value = 433

# Sheet 15: Woodlawn
$ws = $wb.Worksheets.Item('Woodlawn')
$ws.Range('I3').Value = 89
$ws.Range('I7').Value = 248

# Sheet 16: North Lawndale
$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('I3').Value = 181
$ws.Range('I4').Value = 26
$ws.Range('I6').Value = 170
$ws.Range('I7').Value = 512

# Sheet 19: New City
$ws = $wb.Worksheets.Item('New City')
$ws.Range('I6').Value = 94
$ws.Range('I7').Value = 296

# Sheet 20: South Chicago
$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('I3').Value = 106
$ws.Range('I4').Value = 11
$ws.Range('I7').Value = 271

# Sheet 21: West Pullman
$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('I2').Value = 77
$ws.Range('I3').Value = 80
$ws.Range('I7').Value = 213

# Sheet 22: Garfield Park
$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('I3').Value = 218
$ws.Range('I6').Value = 190
$ws.Range('I7').Value = 601

# Sheet 23: Lincoln Park
$ws = $wb.Worksheets.Item('Lincoln Park')
$ws.Range('I3').Value = 12
$ws.Range('I7').Value = 110

# Sheet 24: Loop
$ws = $wb.Worksheets.Item('Loop')
$ws.Range('I2').Value = 70
$ws.Range('I3').Value = 59
$ws.Range('I7').Value = 305

# Sheet 25: Englewood
$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('I2').Value = 255
$ws.Range('I3').Value = 295
$ws.Range('I4').Value = 40
$ws.Range('I6').Value = 238
$ws.Range('I7').Value = 863

# Sheet 26: Chatham
$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('I2').Value = 136
$ws.Range('I3').Value = 104
$ws.Range('I7').Value = 360

# Sheet 27: Irving Park
$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range('I2').Value = 34
$ws.Range('I7').Value = 96

# Sheet 29: River North
$ws = $wb.Worksheets.Item('River North')
$ws.Range('I4').Value = 25
$ws.Range('I7').Value = 202

# Sheet 30: Ashburn
$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('I6').Value = 18
$ws.Range('I7').Value = 91

# Sheet 31: Hermosa
$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range('I3').Value = 20
$ws.Range('I7').Value = 59

# Sheet 32: Humboldt Park
$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range('I2').Value = 127
$ws.Range('I3').Value = 159
$ws.Range('I7').Value = 464

# Sheet 35: Rogers Park
$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('I3').Value = 49
$ws.Range('I7').Value = 193

# Sheet 39: Douglas
$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('I3').Value = 44
$ws.Range('I6').Value = 39
$ws.Range('I7').Value = 128

# Sheet 40: Washington Park
$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range('I6').Value = 52
$ws.Range('I7').Value = 164

# Sheet 41: Chinatown
$ws = $wb.Worksheets.Item('Chinatown')
$ws.Range('I6').Value = 57
$ws.Range('I7').Value = 74

# Sheet 42: Roseland
$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('I2').Value = 107
$ws.Range('I6').Value = 108
$ws.Range('I7').Value = 361

# Sheet 44: Chicago Lawn
$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('I3').Value = 100
$ws.Range('I6').Value = 102
$ws.Range('I7').Value = 326

# Sheet 47: Grand Boulevard
$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('I3').Value = 57
$ws.Range('I4').Value = 8
$ws.Range('I7').Value = 185

# Sheet 51: West Loop
$ws = $wb.Worksheets.Item('West Loop')
$ws.Range('I6').Value = 71
$ws.Range('I7').Value = 127

# Sheet 53: Kenwood
$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('I4').Value = 8
$ws.Range('I7').Value = 94

# Sheet 55: Wicker Park
$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Range('I6').Value = 54
$ws.Range('I7').Value = 88

# Sheet 59: North Center
$ws = $wb.Worksheets.Item('North Center')
$ws.Range('I2').Value = 10
$ws.Range('I7').Value = 36

# Sheet 62: Portage Park
$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('I6').Value = 30
$ws.Range('I7').Value = 111

# Sheet 65: West Town
$ws = $wb.Worksheets.Item('West Town')
$ws.Range('I6').Value = 59
$ws.Range('I7').Value = 101

# Sheet 66: West Elsdon
$ws = $wb.Worksheets.Item('West Elsdon')
$ws.Range('I2').Value = 15
$ws.Range('I6').Value = 16
$ws.Range('I7').Value = 41

# Sheet 72: Streeterville
$ws = $wb.Worksheets.Item('Streeterville')
$ws.Range('I6').Value = 20
$ws.Range('I7').Value = 84

# Sheet 78: Morgan Park
$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('I2').Value = 17
$ws.Range('I7').Value = 66

# Sheet 79: Hyde Park
$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('I6').Value = 67
$ws.Range('I7').Value = 115

# Sheet 82: Old Town
$ws = $wb.Worksheets.Item('Old Town')
$ws.Range('I6').Value = 26
$ws.Range('I7').Value = 49

# Sheet 90: Archer Heights
$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('I6').Value = 15
$ws.Range('I7').Value = 51

# Sheet 91: Beverly
$ws = $wb.Worksheets.Item('Beverly')
$ws.Range('I3').Value = 5
$ws.Range('I7').Value = 27
